$wb = $excel.ActiveWorkbook

# Fix capitalization error: "ProductTb" -> "ProductTB"
$ws = $wb.Worksheets.Item("ProductTb")
$ws.Name = "ProductTB"

# Make the corrected sheet the active/selected tab
$ws.Activate()
